# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Price values are stored as plain text (matches the source feed), so any
# cell whose new text would otherwise be auto-parsed as a number is entered
# with a leading apostrophe (Excel's text quote-prefix) to keep it literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.633.91"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "1.873.36"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'246.76"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "'0.4740"
$ws.Range("E7").Value = "  +0.39%  "

# Row 8
$ws.Range("D8").Value = "'0.2913"
$ws.Range("E8").Value = "  +1.38%  "

# Row 9
$ws.Range("D9").Value = "'0.06485"
$ws.Range("E9").Value = "  +0.24%  "

# Row 10
$ws.Range("D10").Value = "'22.09"
$ws.Range("E10").Value = "  +4.88%  "

# Row 11
$ws.Range("D11").Value = "'0.07726"
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
$ws.Range("D12").Value = "'0.7431"
$ws.Range("E12").Value = "  +4.83%  "

# Row 13
$ws.Range("D13").Value = "'96.77"
$ws.Range("E13").Value = "  +1.97%  "

# Row 14
$ws.Range("D14").Value = "1.871.51"
$ws.Range("E14").Value = "  -2.08%  "

# Row 15
$ws.Range("D15").Value = "'5.150"
$ws.Range("E15").Value = "  +1.03%  "

# Row 16
$ws.Range("D16").Value = "'274.83"
$ws.Range("E16").Value = "  -0.28%  "

# Row 17
$ws.Range("D17").Value = "30.646.41"
$ws.Range("E17").Value = "  +1.26%  "

# Row 18
$ws.Range("D18").Value = "'13.41"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("D20").Value = "'0.000007514"
$ws.Range("E20").Value = "  -0.57%  "

# Row 21
$ws.Range("D21").Value = "2.118.64"

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").Value = "'5.280"
$ws.Range("E23").Value = "  +1.28%  "

# Row 24
$ws.Range("D24").Value = "'6.191"
$ws.Range("E24").Value = "  +0.66%  "

# Row 25
$ws.Range("D25").Value = "'9.261"
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("D26").Value = "'163.57"
$ws.Range("E26").Value = "  -0.94%  "

# Row 27
$ws.Range("D27").Value = "'18.81"
$ws.Range("E27").Value = "  -0.35%  "

# Row 28
$ws.Range("D28").Value = "'1.920"
$ws.Range("E28").Value = "  +0.81%  "

# Row 29
$ws.Range("D29").Value = "'0.09992"
$ws.Range("E29").Value = "  +1.49%  "

# Row 30
$ws.Range("D30").Value = "'1.349"
$ws.Range("E30").Value = "  -1.79%  "

# Row 31
$ws.Range("D31").Value = "'1.508"
$ws.Range("E31").Value = "  -0.44%  "

# Row 32
$ws.Range("D32").Value = "'4.302"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33
$ws.Range("D33").Value = "'4.127"
$ws.Range("E33").Value = "  +2.37%  "

# Row 34
$ws.Range("D34").Value = "'0.04792"
$ws.Range("E34").Value = "  +0.78%  "

# Row 35
$ws.Range("D35").Value = "'1.120"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("D36").Value = "'0.6986"
$ws.Range("E36").Value = "  +0.74%  "

# Row 37
$ws.Range("D37").Value = "'0.9996"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").Value = "'2.711"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("D39").Value = "'0.01854"
$ws.Range("E39").Value = "  +0.65%  "

# Row 40
$ws.Range("D40").Value = "'2.754"
$ws.Range("E40").Value = "  +0.35%  "

# Row 41
$ws.Range("D41").Value = "'6.223"
$ws.Range("E41").Value = "  -1.43%  "

# Row 42
$ws.Range("D42").Value = "'73.63"
$ws.Range("E42").Value = "  +4.41%  "

# Row 43
$ws.Range("D43").Value = "'1.976"
$ws.Range("E43").Value = "  +3.83%  "

# Row 44
$ws.Range("D44").Value = "'0.4185"
$ws.Range("E44").Value = "  +2.00%  "

# Row 45
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("D46").Value = "'0.8345"
$ws.Range("E46").Value = "  -0.76%  "

# Row 47
$ws.Range("D47").Value = "'102.49"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48
$ws.Range("D48").Value = "'9.348"
$ws.Range("E48").Value = "  +1.28%  "

# Row 49
$ws.Range("D49").Value = "'931.16"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'6.974"
$ws.Range("E50").Value = "  -1.63%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.30"
$ws.Range("E51").Value = "  +0.23%  "
